# Finished the data conversion process:
# Add two new columns (C = Density, D = Pressure) that hold the A/B
# datasets converted into geometric units (c = 1, km and solar mass),
# i.e. every value divided by the constant 1.4765679173556.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -------------------------------------------------
$ws.Range("C1").Value = "Density"
$ws.Range("D1").Value = "Pressure"

# --- Row 2: typed in individually (not part of the fill-down) ----------
$ws.Range("C2").Formula = "=A2/1.4765679173556"
$ws.Range("D2").Formula = "=B2/1.4765679173556"

# --- Fill the formula down the rest of the table ------------------------
# First big fill (rows 3-66) done as one rectangular block so Excel keeps
# C and D together in a single shared-formula group.
$ws.Range("C3:D66").Formula = "=A3/1.4765679173556"

# Remaining rows (67-75) filled afterwards -> separate shared-formula group.
$ws.Range("C67:D75").Formula = "=A67/1.4765679173556"

# --- Restore the view: scrolled down near the bottom of the data with
# F72 as the active cell (matches the saved workbook state). -------------
[void]$excel.Goto($ws.Range("A43"), $true)
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F72").Select()
